$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 173 (pushes the former rows 173..287 down to 175..289)
$ws.Rows.Item(173).Insert()
$ws.Rows.Item(173).Insert()

# New row 173: Coliflor "Primera" record dated 2022-08-04 (serial 44777)
$ws.Range("A173").Value = 5
$ws.Range("B173").Value = "Macroferia Regional de Talca"
$ws.Range("C173").Value = "Maule"
$ws.Range("D173").Value = 44777
$ws.Range("E173").Value = 7
$ws.Range("F173").Value = 100112008
$ws.Range("G173").Value = "Coliflor"
$ws.Range("H173").Value = "Sin especificar"
$ws.Range("I173").Value = "Primera"
$ws.Range("J173").Value = 3000
$ws.Range("K173").Value = 1000
$ws.Range("L173").Value = 1000
$ws.Range("M173").Value = 1000
$ws.Range("N173").Value = "$/unidad"
$ws.Range("O173").Value = "Región del Maule"
$ws.Range("P173").Value = 1000
$ws.Range("Q173").Value = 1
$ws.Range("R173").Value = "Hortaliza"

# New row 174: Coliflor "Segunda" record, same date (serial 44777)
$ws.Range("A174").Value = 5
$ws.Range("B174").Value = "Macroferia Regional de Talca"
$ws.Range("C174").Value = "Maule"
$ws.Range("D174").Value = 44777
$ws.Range("E174").Value = 7
$ws.Range("F174").Value = 100112008
$ws.Range("G174").Value = "Coliflor"
$ws.Range("H174").Value = "Sin especificar"
$ws.Range("I174").Value = "Segunda"
$ws.Range("J174").Value = 2000
$ws.Range("K174").Value = 800
$ws.Range("L174").Value = 800
$ws.Range("M174").Value = 800
$ws.Range("N174").Value = "$/unidad"
$ws.Range("O174").Value = "Región del Maule"
$ws.Range("P174").Value = 800
$ws.Range("Q174").Value = 1
$ws.Range("R174").Value = "Hortaliza"
